# Refresh the cryptos price/volume snapshot (GitHub Actions data pull).
# For numeric-looking values in the Price column we prefix with a literal
# apostrophe so Excel stores them as text (preserving formats like
# trailing zeros, e.g. "1.00") instead of silently coercing to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.403.77"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "3.801.45"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'596.39"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "'168.46"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").Value = "3.800.20"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").Value = "'6.40"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").Value = "'0.454"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("D14").Value = "'36.30"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").Value = "4.439.20"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "3.759.11"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "68.464.89"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").Value = "'7.01"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").Value = "'465.40"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E24").Value = "  +7.78%  "
$ws.Range("D25").Value = "'83.99"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("E26").Value = "  -3.60%  "
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").Value = "'10.10"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'30.07"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'7.28"
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("E33").Value = "  -3.68%  "
$ws.Range("D34").Value = "'9.15"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").Value = "3.754.55"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").Value = "'43.86"
$ws.Range("E44").Value = "  +15.10%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.302"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("D46").Value = "'46.91"
$ws.Range("E46").Value = "  +2.94%  "
$ws.Range("D47").Value = "'1.91"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("D49").Value = "'147.06"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").Value = "'390.41"
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("D51").Value = "2.774.44"
$ws.Range("E51").Value = "  +4.22%  "
